$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Where can I find application contact?"
$ws.Range("A8").Select()
